# Update countries & provincias Spain
#
# Applies the 23-Abr-2020 19:52 -> 20:22 data refresh:
#  - Updates the "Datos actualizados..." timestamp banner.
#  - Refreshes the numeric stats for a handful of countries.
#  - Egipto's new total (3891) overtakes Luxemburgo's (3665); Gibraltar's new
#    total (133) overtakes Birmania's (132) -- in both cases the two rows
#    swap positions in the (descending, by "Casos totales") ranking, so the
#    row that used to hold the lower-ranked country keeps that country's old
#    data while the row above it takes on the newly-risen country's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Banner timestamp ------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 20:22"

# --- Francia (row 7) --------------------------------------------------
$ws.Range("B7").Value = 161530
$ws.Range("C7").Value = 1653
$ws.Range("D7").Value = 42088
$ws.Range("E7").Value = 97586
$ws.Range("F7").Value = 5053
$ws.Range("G7").Value = 516
$ws.Range("H7").Value = 21856

# --- Canada (row 16) ------------------------------------------------------
$ws.Range("E16").Value = 25664
$ws.Range("G16").Value = 167
$ws.Range("H16").Value = 2141

# --- Peru (row 21) ---------------------------------------------------------
$ws.Range("B21").Value = 20914
$ws.Range("C21").Value = 1664
$ws.Range("E21").Value = 13315
$ws.Range("G21").Value = 42
$ws.Range("H21").Value = 572

# --- Irlanda (row 23) --------------------------------------------------
$ws.Range("E23").Value = 6644
$ws.Range("G23").Value = 25
$ws.Range("H23").Value = 794

# --- Noruega (row 41) ----------------------------------------------------
$ws.Range("E41").Value = 7136
$ws.Range("F41").Value = 50
$ws.Range("G41").Value = 6
$ws.Range("H41").Value = 193

# --- Egipto / Luxemburgo swap (rows 54-55) ---------------------------------
# Row 54 becomes Egipto (new data), row 55 becomes Luxemburgo (old row-54 data).
$ws.Range("A54").Value = "Egipto"
$ws.Range("B54").Value = 3891
$ws.Range("C54").Value = 232
$ws.Range("D54").Value = 1004
$ws.Range("E54").Value = 2600
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 11
$ws.Range("H54").Value = 287

$ws.Range("A55").Value = "Luxemburgo"
$ws.Range("B55").Value = 3665
$ws.Range("C55").Value = 11
$ws.Range("D55").Value = 728
$ws.Range("E55").Value = 2854
$ws.Range("F55").Value = 27
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 83

# --- Uzbekistan (row 68) --------------------------------------------------
$ws.Range("B68").Value = 1758
$ws.Range("C68").Value = 42
$ws.Range("E68").Value = 1190

# --- Ruanda (row 133) -------------------------------------------------
$ws.Range("B133").Value = 154
$ws.Range("C133").Value = 1
$ws.Range("D133").Value = 87
$ws.Range("E133").Value = 67

# --- Gibraltar / Birmania swap (rows 136-137) -------------------------------
# Row 136 becomes Gibraltar (new data), row 137 becomes Birmania (old row-136 data).
$ws.Range("A136").Value = "Gibraltar"
$ws.Range("B136").Value = 133
$ws.Range("C136").Value = 1
$ws.Range("D136").Value = 129
$ws.Range("E136").Value = 4
$ws.Range("H136").Value = 0

$ws.Range("A137").Value = "Birmania"
$ws.Range("C137").Value = 9
$ws.Range("D137").Value = 9
$ws.Range("E137").Value = 118
$ws.Range("H137").Value = 5
